$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.292.78'
$ws.Range("E2").Value = '  +1.34%  '
$ws.Range("D3").Value = '1.842.69'
$ws.Range("E3").Value = '  +0.59%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9990'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.39'
$ws.Range("E5").Value = '  -0.57%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6873'
$ws.Range("E6").Value = '  -0.92%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9997'
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3032'
$ws.Range("E8").Value = '  -0.24%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07533'
$ws.Range("E9").Value = '  -1.78%  '
$ws.Range("E10").Value = '  +0.14%  '
$ws.Range("E11").Value = '  -1.49%  '
$ws.Range("D12").Value = '1.840.43'
$ws.Range("E12").Value = '  +0.48%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.093'
$ws.Range("E13").Value = '  +0.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6872'
$ws.Range("E14").Value = '  +0.85%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '88.69'
$ws.Range("E15").Value = '  -4.60%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.271'
$ws.Range("E16").Value = '  -3.70%  '
$ws.Range("D17").Value = '29.280.61'
$ws.Range("E17").Value = '  +1.29%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008222'
$ws.Range("D19").Value = '2.088.19'
$ws.Range("E19").Value = '  +0.56%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '232.54'
$ws.Range("E20").Value = '  -3.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9997'
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.484'
$ws.Range("E23").Value = '  +0.49%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9997'
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1465'
$ws.Range("E25").Value = '  -2.14%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '160.04'
$ws.Range("E26").Value = '  +1.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.843'
$ws.Range("E27").Value = '  +1.28%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.13'
$ws.Range("E28").Value = '  -0.25%  '
$ws.Range("E29").Value = '  -1.06%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.273'
$ws.Range("E30").Value = '  +1.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.163'
$ws.Range("E31").Value = '  +0.73%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.205'
$ws.Range("E32").Value = '  +1.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05150'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7735'
$ws.Range("E34").Value = '  +0.05%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.847'
$ws.Range("E35").Value = '  -0.07%  '
$ws.Range("E36").Value = '  -0.02%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.671'
$ws.Range("E37").Value = '  -0.91%  '
$ws.Range("D38").Value = '1.310.29'
$ws.Range("E38").Value = '  +2.57%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01847'
$ws.Range("E39").Value = '  -0.54%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.705'
$ws.Range("E40").Value = '  +0.34%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9437'
$ws.Range("E41").Value = '  -1.27%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.811'
$ws.Range("E42").Value = '  -5.16%  '
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '104.97'
$ws.Range("E43").Value = '  -1.85%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9990'
$ws.Range("E44").Value = '  -0.12%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '65.72'
$ws.Range("E45").Value = '  +3.16%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.682'
$ws.Range("E46").Value = '  +0.13%  '
$ws.Range("B47").Value = 'RocketPoolETH'
$ws.Range("C47").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D47").Value = '1.988.85'
$ws.Range("E47").Value = '  +0.73%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5189'
$ws.Range("E48").Value = '  +0.51%  '
$ws.Range("E49").Value = '  +1.47%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00000000120'
$ws.Range("E50").Value = '  -2.63%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05924'
$ws.Range("E51").Value = '  +0.97%  '
